# Updated BOM and zip
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the "Polarised capacitor" description (Description column, rows 5 and 8)
$ws.Range("G5").Value = "Polarised capacitor "
$ws.Range("G8").Value = "Polarised capacitor "

# Update resistor value and MPN for R5, R6 (row 19): Value column B, MPN column C
$ws.Range("B19").Value = "76K8"
$ws.Range("C19").Value = "CRCW060376K8FKEA"

# Select C9 as the active cell (matches final saved selection state)
$ws.Range("C9").Select()
